# CS297 Time Tracking workbook update
# - Log Wednesday hours (5) for Georgia Fox on the "Week 1" sheet.
#   This single input cascades through the existing SUM-style formulas on
#   "Week 1" (column E running totals) and the "Totals" sheet (which pulls
#   'Week 1'!E$10 and rolls it across the weekly total columns).
# - Leave the active-cell selection on "Week 1" at D8, matching where the
#   user's cursor ended up after entering the value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Week 1")

$ws.Range("D6").Value = 5

$ws.Range("D8").Select() | Out-Null
